$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.573.09"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.340.08"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.14"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.13"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.38"
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.73"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "2.704.14"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "2.352.28"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.807"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "43.465.65"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.14"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.36"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.54"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.06"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.71"
$ws.Range("E29").Value = "  -6.24%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.08"
$ws.Range("E30").Value = "  -8.68%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.28"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.55"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("E35").Value = "  -4.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").Value = "  -4.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.99"
$ws.Range("E37").Value = "  -7.34%  "
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.90"
$ws.Range("E39").Value = "  -7.45%  "
$ws.Range("E40").Value = "  -5.93%  "
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.992.06"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0286"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.59"
$ws.Range("E46").Value = "  -7.04%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.96"
$ws.Range("E47").Value = "  -6.54%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.99"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.48"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.91"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.565.47"
$ws.Range("E51").Value = "  +0.44%  "
